$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the description text in E2 to append " Hi Sara"
$ws.Range("E2").Value = "Beautiful interactions of strings, weaving endlessly into the infinite. Hi Sara"

# Move the active selection to E3 (as in the saved file)
$ws.Range("E3").Select()

# Update the "best fit" column widths to match content change
$ws.Columns.Item(2).ColumnWidth = 12.85546875
$ws.Columns.Item(3).ColumnWidth = 9.85546875
$ws.Columns.Item(5).ColumnWidth = 62.85546875
